# The edit swaps the content of rows 17 and 18 (keeping location/date/etc.
# columns, which were already identical between the two rows, untouched).
# Row 17 becomes the old row 18 data, and row 18 becomes the old row 17 data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17 gets the data that used to be in row 18 ---
$ws.Range("A17").Value = 111961472
$ws.Range("B17").Value = 90857
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 5448
$ws.Range("F17").Value = "Svartvit taggsvamp"
$ws.Range("G17").Value = "Phellodon connatus"
$ws.Range("H17").Value = "(Schultz) nom.prov"
$ws.Range("I17").Value = ""
$ws.Range("J17").Value = ""
$ws.Range("AF17").Value = ""
$ws.Range("AX17").Value = "Stefan Phalagorn Bergström, Annika  Carlberg , Andreas Estensen, Ola Elleström, Anne Järvinen, Emma Sewell, Thomas Strid"

# --- Row 18 gets the data that used to be in row 17 ---
$ws.Range("A18").Value = 111961716
$ws.Range("B18").Value = 81207
$ws.Range("D18").Value = "LC"
$ws.Range("E18").Value = 5046
$ws.Range("F18").Value = "Grön jordtunga"
$ws.Range("G18").Value = "Microglossum viride"
$ws.Range("H18").Value = "(Pers.:Fr.) Gillet"
$ws.Range("I18").NumberFormat = "@"
$ws.Range("I18").Value = "2"
$ws.Range("I18").NumberFormat = "General"
$ws.Range("J18").Value = "mycel"
$ws.Range("AF18").Value = "mikroskoperad"
$ws.Range("AX18").Value = "Stefan Phalagorn Bergström, Andreas Estensen, Annika  Carlberg , Ola Elleström, Thomas Strid, Anne Järvinen, Emma Sewell"
